$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds values as text in the source data (e.g. "0.0000190",
# "5.07", multi-dot big numbers like "76.525.12"). Force text format first so
# Excel does not auto-convert the assigned strings into numbers and mangle them.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "76.474.10"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "3.031.56"
$ws.Range("E3").Value = "  +4.67%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "200.33"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "628.85"
$ws.Range("E6").Value = "  +5.45%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "3.032.20"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "5.06"
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").Value = "3.593.06"
$ws.Range("E14").Value = "  +4.62%  "
$ws.Range("D15").Value = "29.31"
$ws.Range("E15").Value = "  +7.69%  "
$ws.Range("D16").Value = "76.386.59"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "0.0000189"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "3.025.07"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").Value = "13.44"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").Value = "9.06"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").Value = "375.22"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "4.35"
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").Value = "2.26"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("E24").Value = "  +4.72%  "
$ws.Range("D25").Value = "73.13"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "4.38"
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("D28").Value = "9.83"
$ws.Range("E28").Value = "  +2.69%  "
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").Value = "0.995"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "8.33"
$ws.Range("E31").Value = "  +8.65%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "511.64"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +7.59%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "20.78"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "163.98"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").Value = "193.89"
$ws.Range("E38").Value = "  +8.41%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "20.01"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "0.382"
$ws.Range("E40").Value = "  +11.15%  "
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "5.04"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").Value = "42.66"
$ws.Range("E45").Value = "  +6.27%  "
$ws.Range("D46").Value = "1.26"
$ws.Range("E46").Value = "  +5.87%  "
$ws.Range("D47").Value = "1.65"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.711"
$ws.Range("E48").Value = "  +8.46%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.602"
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").Value = "3.88"
$ws.Range("E51").Value = "  +4.70%  "
